$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old sheet had a 2-row header (row1: MW/GWh group labels, row2: Hiver/Ete/Annee
# unit labels). The new layout folds everything into a single header row and adds
# five new leading columns (idx, idx2, Name, Date Start, Date End). Deleting the old
# row 2 shifts all the data rows up by one (old row 3 -> new row 2, etc.) and keeps
# their existing A:E values intact, since those already held idx/idx2/Name/Start/End.
$ws.Rows.Item(2).Delete()

# Rebuild row 1 as a single header row.
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"
$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

# Give the unit-row cells F1:K1 their own (new) font style, matching the rest of
# the workbook's "header" font (Arial 9) used elsewhere for labels.
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Selection moves to the first data row after the edit.
$ws.Range("A2:K2").Select()
